$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.270729541778564
$ws.Range("B1").Value = 2.643998622894287
$ws.Range("C1").Value = 1.869109988212585
$ws.Range("D1").Value = 1.701571106910706
$ws.Range("E1").Value = 1.733711242675781
